{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that holds \"*Patreon\" \u2014 the insertion point per the diff.\nlet patreon = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.trim() === \"*Patreon\") {\n        patreon = paragraphs.items[i];\n        break;\n    }\n}\nif (!patreon) {\n    throw new Error(\"Could not find the '*Patreon' paragraph\");\n}\n\n// Insert two new paragraphs right after it, in order: \"*Kickstarter\", then \"*Indiegogo\".\nconst kickstarter = patreon.insertParagraph(\"*Kickstarter\", Word.InsertLocation.after);\nkickstarter.insertParagraph(\"*Indiegogo\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that contains \"*Patreon\" (the insertion point per the diff).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`a\") -eq \"*Patreon\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the '*Patreon' paragraph\"\n}\n\n# Insert a new paragraph right after \"*Patreon\" and fill it with \"*Kickstarter\".\n$target.Range.InsertParagraphAfter()\n$kickstarter = $target.Next()\n$kickstarter.Range.Text = \"*Kickstarter\"\n\n# Insert another new paragraph right after \"*Kickstarter\" and fill it with \"*Indiegogo\".\n$kickstarter.Range.InsertParagraphAfter()\n$indiegogo = $kickstarter.Next()\n$indiegogo.Range.Text = \"*Indiegogo\"\n"}
